$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.136.57"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "2.992.46"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.31"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.63"
$ws.Range("E6").Value = "  -3.45%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -1.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.31"
$ws.Range("E10").Value = "  -3.19%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0858"
$ws.Range("E12").Value = "  -3.70%  "

$ws.Range("E13").Value = "  -3.30%  "

$ws.Range("D14").Value = "3.459.97"
$ws.Range("E14").Value = "  +1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.71"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").Value = "3.003.58"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").Value = "52.191.91"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.66"
$ws.Range("E21").Value = "  -5.28%  "

$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("E23").Value = "  -2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.73"
$ws.Range("E24").Value = "  -2.31%  "

$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.179"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.93"
$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("E28").Value = "  +2.06%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.110"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.31"
$ws.Range("E31").Value = "  -3.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.26"
$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.19"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.84"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("E36").Value = "  -2.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  -4.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.89"
$ws.Range("E39").Value = "  -4.97%  "

$ws.Range("E40").Value = "  -3.70%  "

$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.72"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.21"
$ws.Range("E44").Value = "  +8.04%  "

$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("D46").Value = "2.126.94"
$ws.Range("E46").Value = "  -1.91%  "

$ws.Range("E47").Value = "  -4.51%  "

$ws.Range("E48").Value = "  -5.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.240"
$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  -3.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.933"
$ws.Range("E51").Value = "  -0.31%  "
